{"js": "// Update the document's page margins (Page Setup), matching the new layout\n// used for the cover page / body after the \"Bluetooth\" content was added:\n//   top:    1440 twips (72 pt)   -> 873 twips (43.65 pt)\n//   left:   1440 twips (72 pt)   -> 2268 twips (113.4 pt)\n//   header: 708  twips (35.4 pt) -> 709  twips (35.45 pt)\n//   footer: 708  twips (35.4 pt) -> 709  twips (35.45 pt)\n// (right/bottom/gutter are unchanged)\n//\n// Office.js expresses PageSetup distances in points (1 pt = 20 twips), so the\n// twips values above are converted to points below.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  // Some hosts expose this as the method `getPageSetup()` (current Word\n  // JS API), others as a plain `pageSetup` property; support both.\n  const pageSetup =\n    typeof section.getPageSetup === \"function\"\n      ? section.getPageSetup()\n      : section.pageSetup;\n  pageSetup.topMargin = 43.65;\n  pageSetup.leftMargin = 113.4;\n  pageSetup.headerDistance = 35.45;\n  pageSetup.footerDistance = 35.45;\n}\n\nawait context.sync();\n", "ps1": "# Update the document's page margins (Page Setup), matching the new layout\n# used for the cover page / body after the \"Bluetooth\" content was added:\n#   top:    1440 twips (72 pt)   -> 873 twips (43.65 pt)\n#   left:   1440 twips (72 pt)   -> 2268 twips (113.4 pt)\n#   header: 708  twips (35.4 pt) -> 709  twips (35.45 pt)\n#   footer: 708  twips (35.4 pt) -> 709  twips (35.45 pt)\n# (right/bottom/gutter are unchanged)\n#\n# Word COM expresses PageSetup distances in points (1 pt = 20 twips), so the\n# twips values above are converted to points below.\n\n$d = $word.ActiveDocument\n\nforeach ($section in $d.Sections) {\n    $ps = $section.PageSetup\n    $ps.TopMargin = 43.65\n    $ps.LeftMargin = 113.4\n    $ps.HeaderDistance = 35.45\n    $ps.FooterDistance = 35.45\n}\n"}
